$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the surviving tail of the original dataset (old rows 2-13)
# before we overwrite the sheet, since those rows end up shifted down to
# rows 20-31 in the final layout (old rows 14-21 are dropped).
$keepA = @()
$keepB = @()
$keepC = @()
for ($r = 2; $r -le 13; $r++) {
    $keepA += ,($ws.Range("A$r").Value2)
    $keepB += ,($ws.Range("B$r").Value2)
    $keepC += ,($ws.Range("C$r").Value2)
}

# Write the new accelerometer samples (captured May 9th) into rows 2-19.
$ws.Range("A2").Value = -3.58582592010498
$ws.Range("B2").Value = 8.237812042236328
$ws.Range("C2").Value = 0.2284512519836425
$ws.Range("A3").Value = -3.345468521118164
$ws.Range("B3").Value = 8.276437759399414
$ws.Range("C3").Value = -0.0355764925479888
$ws.Range("A4").Value = -3.182104587554932
$ws.Range("B4").Value = 8.366311073303223
$ws.Range("C4").Value = 0.2440185844898224
$ws.Range("A5").Value = -3.532341480255127
$ws.Range("B5").Value = 8.214254379272461
$ws.Range("C5").Value = 0.022090196609497
$ws.Range("A6").Value = -3.395359516143799
$ws.Range("B6").Value = 8.257164001464844
$ws.Range("C6").Value = 0.0318913161754608
$ws.Range("A7").Value = -3.210868835449219
$ws.Range("B7").Value = 8.186161994934082
$ws.Range("C7").Value = -0.0299544632434844
$ws.Range("A8").Value = -3.197915077209473
$ws.Range("B8").Value = 8.2127103805542
$ws.Range("C8").Value = 0.0095521211624145
$ws.Range("A9").Value = -3.255885124206543
$ws.Range("B9").Value = 8.21357250213623
$ws.Range("C9").Value = -0.070087194442749
$ws.Range("A10").Value = -3.168186187744141
$ws.Range("B10").Value = 8.229218482971191
$ws.Range("C10").Value = -0.0939638018608093
$ws.Range("A11").Value = -3.290424346923828
$ws.Range("B11").Value = 8.159111976623535
$ws.Range("C11").Value = -0.1956053972244262
$ws.Range("A12").Value = -3.586246490478516
$ws.Range("B12").Value = 8.05996036529541
$ws.Range("C12").Value = -0.0569053888320922
$ws.Range("A13").Value = -3.300580024719238
$ws.Range("B13").Value = 8.124805450439453
$ws.Range("C13").Value = -0.1510338187217712
$ws.Range("A14").Value = -2.94456958770752
$ws.Range("B14").Value = 8.213338851928711
$ws.Range("C14").Value = -0.0691232085227966
$ws.Range("A15").Value = -3.094478130340576
$ws.Range("B15").Value = 8.185011863708496
$ws.Range("C15").Value = -0.1414701342582702
$ws.Range("A16").Value = -3.334782123565674
$ws.Range("B16").Value = 8.0909423828125
$ws.Range("C16").Value = -0.0672928094863891
$ws.Range("A17").Value = -3.368669509887696
$ws.Range("B17").Value = 8.071453094482422
$ws.Range("C17").Value = 0.1171565353870391
$ws.Range("A18").Value = -3.412579536437988
$ws.Range("B18").Value = 8.035589218139648
$ws.Range("C18").Value = 0.0592367351055145
$ws.Range("A19").Value = -3.019937038421631
$ws.Range("B19").Value = 8.054733276367188
$ws.Range("C19").Value = -0.1420263051986694

# Re-write the preserved tail into rows 20-31.
for ($i = 0; $i -lt $keepA.Count; $i++) {
    $r = 20 + $i
    $ws.Range("A$r").Value = $keepA[$i]
    $ws.Range("B$r").Value = $keepB[$i]
    $ws.Range("C$r").Value = $keepC[$i]
}

# The sheet previously ended at row 21; now it ends at row 31, so no
# leftover rows need clearing.
